$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map old "Room <Name>" references to new "Room <Number>" references
# (Room Shelley) -> (Room G22)
# (Room Stephane) -> (Room G19)
# (Room Ivy) -> (Room G14)

$ws.Range("B7").Value = "Rehearsal with pianist`n(Room G22)"
$ws.Range("D7").Value = "Private Lesson with Stephane RETY `n(Room G19)"
$ws.Range("F7").Value = "Flute MasterClass`n(Room G19)"

$ws.Range("C11").Value = "Private Lesson with Stephane RETY & pianist `n(Room G19)"

$ws.Range("B19").Value = "Private Lesson with Ivy CHUANG `n(Room G14)"
$ws.Range("E19").Value = "Private Lesson with Ivy CHUANG `n(Room G14)"
$ws.Range("F19").Value = "Flute MasterClass`n(Room G19)"

$ws.Range("B23").Value = "Ensemble `n(Room G19)"
$ws.Range("C23").Value = "Ensemble `n(Room G19)"
$ws.Range("D23").Value = "Ensemble `n(Room G19)"
$ws.Range("E23").Value = "Ensemble `n(Room G19)"
$ws.Range("F23").Value = "Ensemble `n(Room G19)"
